$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.312.07"
$ws.Range("E2").Value = "'  +2.15%  "
$ws.Range("D3").Value = "'2.547.37"
$ws.Range("E3").Value = "'  +4.69%  "
$ws.Range("E4").Value = "'  -0.08%  "
$ws.Range("D5").Value = "'571.47"
$ws.Range("E5").Value = "'  +3.02%  "
$ws.Range("D6").Value = "'151.08"
$ws.Range("E6").Value = "'  +8.97%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "'  -0.08%  "
$ws.Range("D8").Value = "'0.589"
$ws.Range("E8").Value = "'  +0.97%  "
$ws.Range("D9").Value = "'2.544.60"
$ws.Range("E9").Value = "'  +4.65%  "
$ws.Range("E10").Value = "'  +1.88%  "
$ws.Range("E11").Value = "'  -0.20%  "
$ws.Range("E12").Value = "'  +1.25%  "
$ws.Range("D13").Value = "'0.359"
$ws.Range("E13").Value = "'  +3.26%  "
$ws.Range("D14").Value = "'28.50"
$ws.Range("E14").Value = "'  +9.44%  "
$ws.Range("D15").Value = "'2.998.72"
$ws.Range("E15").Value = "'  +4.58%  "
$ws.Range("D16").Value = "'63.247.78"
$ws.Range("E16").Value = "'  +2.17%  "
$ws.Range("D17").Value = "'0.0000144"
$ws.Range("E17").Value = "'  +0.96%  "
$ws.Range("D18").Value = "'2.546.57"
$ws.Range("E18").Value = "'  +4.44%  "
$ws.Range("D19").Value = "'11.65"
$ws.Range("E19").Value = "'  +4.70%  "
$ws.Range("D20").Value = "'340.50"
$ws.Range("E20").Value = "'  -1.04%  "
$ws.Range("D21").Value = "'4.37"
$ws.Range("E21").Value = "'  +3.98%  "
$ws.Range("D22").Value = "'6.88"
$ws.Range("E22").Value = "'  +0.87%  "
$ws.Range("E23").Value = "'  +0.10%  "
$ws.Range("D24").Value = "'66.01"
$ws.Range("E24").Value = "'  +1.41%  "
$ws.Range("E25").Value = "'  -1.16%  "
$ws.Range("D26").Value = "'1.61"
$ws.Range("E26").Value = "'  +6.98%  "
$ws.Range("E27").Value = "'  +0.21%  "
$ws.Range("D28").Value = "'8.48"
$ws.Range("E28").Value = "'  +3.17%  "
$ws.Range("D29").Value = "'1.48"
$ws.Range("E29").Value = "'  +10.72%  "
$ws.Range("D30").Value = "'7.07"
$ws.Range("E30").Value = "'  +11.92%  "
$ws.Range("D31").Value = "'0.0₃0836"
$ws.Range("E31").Value = "'  +6.38%  "
$ws.Range("D32").Value = "'1.87"
$ws.Range("E32").Value = "'  +3.51%  "
$ws.Range("D33").Value = "'176.69"
$ws.Range("E33").Value = "'  +2.72%  "
$ws.Range("D34").Value = "'1.56"
$ws.Range("E34").Value = "'  +7.59%  "
$ws.Range("D35").Value = "'420.19"
$ws.Range("E35").Value = "'  +14.83%  "
$ws.Range("D36").Value = "'0.407"
$ws.Range("E36").Value = "'  +2.74%  "
$ws.Range("D37").Value = "'19.10"
$ws.Range("E37").Value = "'  +2.84%  "
$ws.Range("D38").Value = "'4.44"
$ws.Range("E38").Value = "'  -1.40%  "
$ws.Range("D40").Value = "'1.76"
$ws.Range("E40").Value = "'  +2.91%  "
$ws.Range("E41").Value = "'  -0.02%  "
$ws.Range("D42").Value = "'39.98"
$ws.Range("E42").Value = "'  +2.10%  "
$ws.Range("D43").Value = "'156.09"
$ws.Range("E43").Value = "'  +6.86%  "
$ws.Range("D44").Value = "'3.81"
$ws.Range("E44").Value = "'  +3.72%  "
$ws.Range("D45").Value = "'21.15"
$ws.Range("E45").Value = "'  +2.56%  "
$ws.Range("D46").Value = "'0.610"
$ws.Range("E46").Value = "'  +3.61%  "
$ws.Range("D47").Value = "'0.0531"
$ws.Range("E47").Value = "'  +2.82%  "
$ws.Range("D48").Value = "'0.0966"
$ws.Range("E48").Value = "'  +1.01%  "
$ws.Range("D49").Value = "'0.0237"
$ws.Range("E49").Value = "'  +7.40%  "
$ws.Range("D50").Value = "'18.60"
$ws.Range("E50").Value = "'  +4.11%  "
$ws.Range("B51").Value = "'dogwifhat"
$ws.Range("C51").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "'1.82"
$ws.Range("E51").Value = "'  +6.76%  "
